$d = $word.ActiveDocument

$ldq = [char]0x201C   # left double quotation mark "
$rdq = [char]0x201D   # right double quotation mark "

# 1. Remove the sentence about not needing to fill out the Affix Transfer Tax
#    Stamp box "until later" (it trailed the "Read through your forms..."
#    bullet point).
$d.Content.Find.Execute(
    " You do not need to fill out the " + $ldq + "Affix Transfer Tax Stamp" + $rdq + " box until later.s",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

# 2. Append new guidance about signing the Affix Transfer Tax Stamp box to the
#    end of the "...must sign and date the TODI..." bullet point.
$d.Content.Find.Execute(
    "you must bring those people with you to be witnesses.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "you must bring those people with you to be witnesses. You should also sign the box labeled " + $ldq + "Affix Transfer Tax Stamp." + $rdq,
    2) | Out-Null
